$wb = $excel.ActiveWorkbook

# --- 1. "Colors" sheet: insert a new top row with a "Columns colors" header,
#        and rename the old "Costs sheet" header to "Costs sheet colors".
$colors = $wb.Worksheets.Item("Colors")

# Inserting a whole row at the top pushes every existing row (and its own
# per-row style) down by one, which is exactly what is needed here: the
# style indices stay attached to their original (now shifted) rows.
$colors.Range("A1").EntireRow.Insert()
$colors.Range("A1").Value = "Columns colors"

# The row that used to hold "Costs sheet" (old A6) is now A7; rename it.
$colors.Range("A7").Value = "Costs sheet colors"

# Column A got one character narrower (50.71 -> 49.71).
$colors.Columns(1).ColumnWidth = 48.83

# --- 2. Refresh the "Created:" timestamp shown on the Costs sheets.
$costs = $wb.Worksheets.Item("Costs")
$costsDnf = $wb.Worksheets.Item("Costs (DNF)")
$costs.Range("B25").Value = "2024-04-23 15:58:50"
$costsDnf.Range("B21").Value = "2024-04-23 15:58:50"
